# Auto-generated edit script: applies the scraped-data update
# (new scrape pass at 11:38:09, additional rows, and a handful of
# reordered/updated rows) to all three worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')

$ws.Cells.Item(2, 1).Value = 'Última actualización: 11:38:09'

$ws.Cells.Item(3, 1).Value = 'Total filas: 132'

$ws.Cells.Item(55, 1).Value = '08:49:06'
$ws.Cells.Item(55, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(55, 4).Value = 4

$ws.Cells.Item(56, 1).Value = '08:14:55'
$ws.Cells.Item(56, 3).Value = '215B_EL PATO'
$ws.Cells.Item(56, 4).Value = 39

$ws.Cells.Item(71, 1).Value = '08:14:55'
$ws.Cells.Item(71, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(71, 4).Value = 77

$ws.Cells.Item(72, 1).Value = '08:49:06'
$ws.Cells.Item(72, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(72, 4).Value = 42

$ws.Cells.Item(103, 1).Value = '11:01:19'
$ws.Cells.Item(103, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(103, 4).Value = 14

$ws.Cells.Item(104, 1).Value = '10:32:07'
$ws.Cells.Item(104, 3).Value = '14_ABASTO'
$ws.Cells.Item(104, 4).Value = 43

$ws.Cells.Item(111, 1).Value = '11:38:09'
$ws.Cells.Item(111, 4).Value = 4

$ws.Cells.Item(112, 1).Value = '11:38:09'
$ws.Cells.Item(112, 4).Value = 7

$ws.Cells.Item(115, 1).Value = '11:38:09'
$ws.Cells.Item(115, 2).Value = '11:52'
$ws.Cells.Item(115, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(115, 4).Value = 14

$ws.Cells.Item(116, 1).Value = '11:38:09'
$ws.Cells.Item(116, 2).Value = '11:53'
$ws.Cells.Item(116, 3).Value = '225_GOMEZ'
$ws.Cells.Item(116, 4).Value = 15

$ws.Cells.Item(117, 1).Value = '11:38:09'
$ws.Cells.Item(117, 2).Value = '11:58'
$ws.Cells.Item(117, 3).Value = '17_ROMERO'
$ws.Cells.Item(117, 4).Value = 20

$ws.Cells.Item(118, 1).Value = '11:01:19'
$ws.Cells.Item(118, 2).Value = '12:05'
$ws.Cells.Item(118, 4).Value = 64

$ws.Cells.Item(119, 1).Value = '11:38:09'
$ws.Cells.Item(119, 2).Value = '12:06'
$ws.Cells.Item(119, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(119, 4).Value = 28

$ws.Cells.Item(120, 1).Value = '11:38:09'
$ws.Cells.Item(120, 3).Value = '15_ABASTO'
$ws.Cells.Item(120, 4).Value = 32

$ws.Cells.Item(121, 1).Value = '11:38:09'
$ws.Cells.Item(121, 2).Value = '12:10'
$ws.Cells.Item(121, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(121, 4).Value = 32

$ws.Cells.Item(122, 1).Value = '11:38:09'
$ws.Cells.Item(122, 2).Value = '12:17'
$ws.Cells.Item(122, 3).Value = '10_OLMOS'
$ws.Cells.Item(122, 4).Value = 39

$ws.Cells.Item(123, 1).Value = '11:38:09'
$ws.Cells.Item(123, 2).Value = '12:22'
$ws.Cells.Item(123, 3).Value = '215C_EL PATO'
$ws.Cells.Item(123, 4).Value = 44

$ws.Cells.Item(124, 1).Value = '11:38:09'
$ws.Cells.Item(124, 2).Value = '12:27'
$ws.Cells.Item(124, 3).Value = '14_ABASTO'
$ws.Cells.Item(124, 4).Value = 49

$ws.Cells.Item(125, 1).Value = '11:38:09'
$ws.Cells.Item(125, 2).Value = '12:31'
$ws.Cells.Item(125, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(125, 4).Value = 53

$ws.Cells.Item(126, 2).Value = '12:32'
$ws.Cells.Item(126, 3).Value = '14_ABASTO'
$ws.Cells.Item(126, 4).Value = 91

$ws.Cells.Item(127, 1).Value = '11:38:09'
$ws.Cells.Item(127, 2).Value = '12:33'
$ws.Cells.Item(127, 3).Value = '14_ABASTO'
$ws.Cells.Item(127, 4).Value = 55
$ws.Cells.Item(127, 5).Value = 'LP1912'

$ws.Cells.Item(128, 1).Value = '11:38:09'
$ws.Cells.Item(128, 2).Value = '12:34'
$ws.Cells.Item(128, 3).Value = '15_ABASTO'
$ws.Cells.Item(128, 4).Value = 56
$ws.Cells.Item(128, 5).Value = 'LP1912'

$ws.Cells.Item(129, 1).Value = '11:38:09'
$ws.Cells.Item(129, 2).Value = '12:37'
$ws.Cells.Item(129, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(129, 4).Value = 59
$ws.Cells.Item(129, 5).Value = 'LP1912'

$ws.Cells.Item(130, 1).Value = '11:38:09'
$ws.Cells.Item(130, 2).Value = '12:48'
$ws.Cells.Item(130, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(130, 4).Value = 70
$ws.Cells.Item(130, 5).Value = 'LP1912'

$ws.Cells.Item(131, 1).Value = '11:38:09'
$ws.Cells.Item(131, 2).Value = '12:48'
$ws.Cells.Item(131, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(131, 4).Value = 70
$ws.Cells.Item(131, 5).Value = 'LP1912'

$ws.Cells.Item(132, 1).Value = '11:38:09'
$ws.Cells.Item(132, 2).Value = '13:03'
$ws.Cells.Item(132, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(132, 4).Value = 85
$ws.Cells.Item(132, 5).Value = 'LP1912'

$ws.Cells.Item(133, 1).Value = '11:38:09'
$ws.Cells.Item(133, 2).Value = '13:04'
$ws.Cells.Item(133, 3).Value = '215C_EL PATO'
$ws.Cells.Item(133, 4).Value = 86
$ws.Cells.Item(133, 5).Value = 'LP1912'

$ws.Cells.Item(134, 1).Value = '11:38:09'
$ws.Cells.Item(134, 2).Value = '13:13'
$ws.Cells.Item(134, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(134, 4).Value = 95
$ws.Cells.Item(134, 5).Value = 'LP1912'

$ws.Cells.Item(135, 1).Value = '11:38:09'
$ws.Cells.Item(135, 2).Value = '13:17'
$ws.Cells.Item(135, 3).Value = '10_OLMOS'
$ws.Cells.Item(135, 4).Value = 99
$ws.Cells.Item(135, 5).Value = 'LP1912'

$ws.Cells.Item(136, 1).Value = '11:38:09'
$ws.Cells.Item(136, 2).Value = '13:25'
$ws.Cells.Item(136, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(136, 4).Value = 107
$ws.Cells.Item(136, 5).Value = 'LP1912'

$ws.Cells.Item(137, 1).Value = '11:38:09'
$ws.Cells.Item(137, 2).Value = '13:33'
$ws.Cells.Item(137, 3).Value = '215A_EL PATO'
$ws.Cells.Item(137, 4).Value = 115
$ws.Cells.Item(137, 5).Value = 'LP1912'


# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')

$ws.Cells.Item(2, 1).Value = 'Última actualización: 11:38:09'

$ws.Cells.Item(3, 1).Value = 'Total filas: 23'

$ws.Cells.Item(25, 1).Value = '11:38:09'
$ws.Cells.Item(25, 4).Value = 4

$ws.Cells.Item(26, 1).Value = '11:38:09'
$ws.Cells.Item(26, 4).Value = 44

$ws.Cells.Item(27, 1).Value = '11:38:09'
$ws.Cells.Item(27, 2).Value = '13:04'
$ws.Cells.Item(27, 3).Value = '215C_EL PATO'
$ws.Cells.Item(27, 4).Value = 86
$ws.Cells.Item(27, 5).Value = 'LP1912'

$ws.Cells.Item(28, 1).Value = '11:38:09'
$ws.Cells.Item(28, 2).Value = '13:33'
$ws.Cells.Item(28, 3).Value = '215A_EL PATO'
$ws.Cells.Item(28, 4).Value = 115
$ws.Cells.Item(28, 5).Value = 'LP1912'


# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')

$ws.Cells.Item(2, 1).Value = 'Última actualización: 11:38:09'

$ws.Cells.Item(3, 1).Value = 'Total filas: 22'

$ws.Cells.Item(19, 1).Value = '08:49:06'
$ws.Cells.Item(19, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(19, 4).Value = 101

$ws.Cells.Item(20, 1).Value = '09:42:42'
$ws.Cells.Item(20, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(20, 4).Value = 48

$ws.Cells.Item(26, 1).Value = '11:38:09'
$ws.Cells.Item(26, 2).Value = '13:12'
$ws.Cells.Item(26, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(26, 4).Value = 94
$ws.Cells.Item(26, 5).Value = 'L6203'

$ws.Cells.Item(27, 1).Value = '11:38:09'
$ws.Cells.Item(27, 2).Value = '13:21'
$ws.Cells.Item(27, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(27, 4).Value = 103
$ws.Cells.Item(27, 5).Value = 'L6173'

